$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# Add "lesson 1" localization entries (rows 26-33) in the same cell-write
# order the original authoring session used, so new shared strings are
# interned in the same sequence.
$ws.Range("A28").Value = "lesson1_distribute_1"
$ws.Range("B26").Value = "Multiplying a number of two or more digits can be tricky, but there’s a special trick to make it easier!"
$ws.Range("A26").Value = "lesson1_intro_1"
$ws.Range("A27").Value = "lesson1_intro_2"
$ws.Range("B27").Value = "Let’s take a quick look at a math concept that we will be using for this trick."
$ws.Range("A29").Value = "lesson1_area_1"
$ws.Range("B29").Value = "To help visualize this, we will consider these smaller pieces as chunks of a rectangle’s area."
$ws.Range("A30").Value = "lesson1_area_2"
$ws.Range("A31").Value = "lesson1_area_3"
$ws.Range("A32").Value = "lesson1_area_4"
$ws.Range("B30").Value = "As you can see, the rectangle is split into two."
$ws.Range("B31").Value = "Then we compute the area of those two rectangles."
$ws.Range("B32").Value = "And finally, we add these two areas together to get the area of the whole rectangle."
$ws.Range("B33").Value = "Let’s give it a try. Connect these two blobs to initiate the attack!"
$ws.Range("A33").Value = "lesson1_connect_1"
$ws.Range("B28").Value = "By using the distributive property, we can split up the large number of an equation into smaller pieces."

$ws.Range("B28").Select()
